$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data to match the latest scrape.
# Column D holds price strings that can look numeric (e.g. "1.100"); force
# the cell's number format to Text ("@") before assigning so Excel keeps
# them as literal strings instead of silently converting to floating point
# numbers (which would drop significant trailing/leading characters).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.653.09'
$ws.Range('E2').Value = '  -2.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.761.02'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.63'
$ws.Range('E5').Value = '  -1.16%  '
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4297'
$ws.Range('E7').Value = '  -1.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3605'
$ws.Range('E8').Value = '  -1.78%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07559'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.19'
$ws.Range('E10').Value = '  -6.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.111'
$ws.Range('E11').Value = '  -2.77%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.77'
$ws.Range('E13').Value = '  -6.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.072'
$ws.Range('E14').Value = '  -3.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.251'
$ws.Range('E15').Value = '  -3.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.759.47'
$ws.Range('E16').Value = '  -4.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.58'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001066'
$ws.Range('E18').Value = '  -1.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06437'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.12'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.881'
$ws.Range('E22').Value = '  -6.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.695.67'
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  -3.47%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.081'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.14'
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.55'
$ws.Range('E27').Value = '  -1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.960.21'
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.159'
$ws.Range('E29').Value = '  -6.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.84'
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.100'
$ws.Range('E31').Value = '  -9.92%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.674'
$ws.Range('E32').Value = '  +5.73%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.606'
$ws.Range('E33').Value = '  -6.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08966'
$ws.Range('E34').Value = '  -2.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.19'
$ws.Range('E35').Value = '  -6.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02299'
$ws.Range('E36').Value = '  -2.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2115'
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06017'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6359'
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.963'
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('E41').Value = '  -0.42%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.397'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.893'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.38'
$ws.Range('E45').Value = '  -4.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5943'
$ws.Range('E46').Value = '  -2.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.714'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.51'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.176'
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06877'
$ws.Range('E51').Value = '  -1.94%  '

Write-Host "Updated cryptos list"
